# ---------------------------------------------------------------------------
# Scheduled-runner market-data refresh for Sheets workbook.
#
# Pulls fresh Universalis-style market board snapshots into the per-job
# columns (H:N => currentAveragePrice / *NQ / *HQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ) on each job sheet. Only the market-derived
# value columns change; leve metadata (A:G) is left untouched.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1031.6072
$ws.Range("I19").Value = 1104.65
$ws.Range("J19").Value = 849
$ws.Range("K19").Value = 1104.65
$ws.Range("L19").Value = 849
$ws.Range("M19").Value = -929.6500000000001
$ws.Range("N19").Value = -1199
$ws.Range("H40").Value = 2025.125
$ws.Range("I40").Value = 2322.4443
$ws.Range("J40").Value = 1642.8572
$ws.Range("K40").Value = 2322.4443
$ws.Range("L40").Value = 1642.8572
$ws.Range("M40").Value = -2147.4443
$ws.Range("N40").Value = -1992.8572
$ws.Range("H99").Value = 857.5333000000001
$ws.Range("I99").Value = 544
$ws.Range("J99").Value = 1719.75
$ws.Range("K99").Value = 1632
$ws.Range("L99").Value = 5159.25
$ws.Range("M99").Value = -134
$ws.Range("N99").Value = -8155.25
$ws.Range("H106").Value = 8203.333000000001
$ws.Range("I106").Value = 9570
$ws.Range("J106").Value = 7930
$ws.Range("K106").Value = 9570
$ws.Range("L106").Value = 7930
$ws.Range("M106").Value = -8939
$ws.Range("N106").Value = -9192
$ws.Range("H116").Value = 4568.6855
$ws.Range("I116").Value = 2890.4
$ws.Range("J116").Value = 5240
$ws.Range("K116").Value = 2890.4
$ws.Range("L116").Value = 5240
$ws.Range("M116").Value = 551.5999999999999
$ws.Range("N116").Value = -12124
$ws.Range("H132").Value = 1646298.1
$ws.Range("I132").Value = 1894904
$ws.Range("J132").Value = 5499.6
$ws.Range("K132").Value = 5684712
$ws.Range("L132").Value = 16498.8
$ws.Range("M132").Value = -5682182
$ws.Range("N132").Value = -21558.8
$ws.Range("H141").Value = 2618.7646
$ws.Range("I141").Value = 926.425
$ws.Range("J141").Value = 8772.727999999999
$ws.Range("K141").Value = 2779.275
$ws.Range("L141").Value = 26318.184
$ws.Range("M141").Value = 2400.725
$ws.Range("N141").Value = -36678.18399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10937.761
$ws.Range("I32").Value = 8089.2856
$ws.Range("J32").Value = 62780
$ws.Range("K32").Value = 8089.2856
$ws.Range("L32").Value = 62780
$ws.Range("M32").Value = -7802.2856
$ws.Range("N32").Value = -63354
$ws.Range("H45").Value = 25005948
$ws.Range("I45").Value = 38469988
$ws.Range("J45").Value = 1302
$ws.Range("K45").Value = 38469988
$ws.Range("L45").Value = 1302
$ws.Range("M45").Value = -38469611
$ws.Range("N45").Value = -2056
$ws.Range("H74").Value = 7354336.5
$ws.Range("I74").Value = 10417912
$ws.Range("J74").Value = 1755.65
$ws.Range("K74").Value = 10417912
$ws.Range("L74").Value = 1755.65
$ws.Range("M74").Value = -10417038
$ws.Range("N74").Value = -3503.65
$ws.Range("H77").Value = 7354336.5
$ws.Range("I77").Value = 10417912
$ws.Range("J77").Value = 1755.65
$ws.Range("K77").Value = 52089560
$ws.Range("L77").Value = 8778.25
$ws.Range("M77").Value = -52085192
$ws.Range("N77").Value = -17514.25
$ws.Range("H132").Value = 2336.8604
$ws.Range("I132").Value = 1603.7858
$ws.Range("J132").Value = 3705.2666
$ws.Range("K132").Value = 4811.357400000001
$ws.Range("L132").Value = 11115.7998
$ws.Range("M132").Value = -2281.357400000001
$ws.Range("N132").Value = -16175.7998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 74900
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 74900
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 74900
$ws.Range("N51").Value = -75882
$ws.Range("H63").Value = 17500
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 17500
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 17500
$ws.Range("N63").Value = -18872
$ws.Range("H66").Value = 17500
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 17500
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 52500
$ws.Range("N66").Value = -59364
$ws.Range("H134").Value = 8622922
$ws.Range("I134").Value = 13159524
$ws.Range("J134").Value = 3378
$ws.Range("K134").Value = 39478572
$ws.Range("L134").Value = 10134
$ws.Range("M134").Value = -39476037
$ws.Range("N134").Value = -15204

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3706.889
$ws.Range("I16").Value = 3030.1428
$ws.Range("J16").Value = 4137.5454
$ws.Range("K16").Value = 3030.1428
$ws.Range("L16").Value = 4137.5454
$ws.Range("M16").Value = -2743.1428
$ws.Range("N16").Value = -4711.5454
$ws.Range("H31").Value = 7815526
$ws.Range("I31").Value = 13514576
$ws.Range("J31").Value = 5716.8887
$ws.Range("K31").Value = 13514576
$ws.Range("L31").Value = 5716.8887
$ws.Range("M31").Value = -13514281
$ws.Range("N31").Value = -6306.8887
$ws.Range("H34").Value = 7815526
$ws.Range("I34").Value = 13514576
$ws.Range("J34").Value = 5716.8887
$ws.Range("K34").Value = 13514576
$ws.Range("L34").Value = 5716.8887
$ws.Range("M34").Value = -13514374
$ws.Range("N34").Value = -6120.8887
$ws.Range("H58").Value = 2020.8368
$ws.Range("I58").Value = 1379.5
$ws.Range("J58").Value = 2745.8262
$ws.Range("K58").Value = 1379.5
$ws.Range("L58").Value = 2745.8262
$ws.Range("M58").Value = -1176.5
$ws.Range("N58").Value = -3151.8262
$ws.Range("H113").Value = 3706.889
$ws.Range("I113").Value = 3030.1428
$ws.Range("J113").Value = 4137.5454
$ws.Range("K113").Value = 3030.1428
$ws.Range("L113").Value = 4137.5454
$ws.Range("M113").Value = -860.1428000000001
$ws.Range("N113").Value = -8477.545399999999
$ws.Range("H132").Value = 2262.1177
$ws.Range("I132").Value = 1622.9474
$ws.Range("J132").Value = 3071.7334
$ws.Range("K132").Value = 4868.8422
$ws.Range("L132").Value = 9215.200199999999
$ws.Range("M132").Value = -2338.8422
$ws.Range("N132").Value = -14275.2002
$ws.Range("H134").Value = 1872.8
$ws.Range("I134").Value = 1829.8462
$ws.Range("J134").Value = 1919.3334
$ws.Range("K134").Value = 5489.5386
$ws.Range("L134").Value = 5758.0002
$ws.Range("M134").Value = -2954.5386
$ws.Range("N134").Value = -10828.0002
$ws.Range("H136").Value = 2020.8368
$ws.Range("I136").Value = 1379.5
$ws.Range("J136").Value = 2745.8262
$ws.Range("K136").Value = 4138.5
$ws.Range("L136").Value = 8237.4786
$ws.Range("M136").Value = -1588.5
$ws.Range("N136").Value = -13337.4786

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 3153.1785
$ws.Range("I139").Value = 1867.8422
$ws.Range("J139").Value = 5866.6665
$ws.Range("K139").Value = 5603.5266
$ws.Range("L139").Value = 17599.9995
$ws.Range("M139").Value = -463.5266000000001
$ws.Range("N139").Value = -27879.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H70").Value = 5615.5713
$ws.Range("I70").Value = 4825
$ws.Range("J70").Value = 6669.6665
$ws.Range("K70").Value = 4825
$ws.Range("L70").Value = 6669.6665
$ws.Range("M70").Value = -4555
$ws.Range("N70").Value = -7209.6665
$ws.Range("H73").Value = 5615.5713
$ws.Range("I73").Value = 4825
$ws.Range("J73").Value = 6669.6665
$ws.Range("K73").Value = 4825
$ws.Range("L73").Value = 6669.6665
$ws.Range("M73").Value = -3889
$ws.Range("N73").Value = -8541.666499999999
$ws.Range("H132").Value = 5966.8438
$ws.Range("I132").Value = 9602.77
$ws.Range("J132").Value = 3479.1052
$ws.Range("K132").Value = 28808.31
$ws.Range("L132").Value = 10437.3156
$ws.Range("M132").Value = -26278.31
$ws.Range("N132").Value = -15497.3156

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3332.3242
$ws.Range("I132").Value = 2658.3157
$ws.Range("J132").Value = 4043.7778
$ws.Range("K132").Value = 7974.9471
$ws.Range("L132").Value = 12131.3334
$ws.Range("M132").Value = -5444.9471
$ws.Range("N132").Value = -17191.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H95").Value = 34883.332
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 34883.332
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 34883.332
$ws.Range("N95").Value = -40375.332
$ws.Range("H132").Value = 1368.2549
$ws.Range("I132").Value = 1019.05884
$ws.Range("J132").Value = 2066.647
$ws.Range("K132").Value = 3057.17652
$ws.Range("L132").Value = 6199.941
$ws.Range("M132").Value = -527.17652
$ws.Range("N132").Value = -11259.941
$ws.Range("H136").Value = 4667.3423
$ws.Range("I136").Value = 1507.92
$ws.Range("J136").Value = 10743.154
$ws.Range("K136").Value = 4523.76
$ws.Range("L136").Value = 32229.462
$ws.Range("M136").Value = -1973.76
$ws.Range("N136").Value = -37329.462
